$wb = $excel.ActiveWorkbook

# The "Replace Substrings" sheet is our template: same Action/Time/Content
# layout, same styles. Duplicate it to create the new "Text Case" sheet,
# placed right after it (matches sheetId=12 / rId12 ordering in the diff).
$src = $wb.Worksheets.Item("Replace Substrings")
$src.Copy([System.Reflection.Missing]::Value, $src)

$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "Text Case"

# Fill in the GOMS content for the "Text Case" task.
$new.Range("A3").Value = "Identify Text Columns"
$new.Range("B3").Value = "2 min"
$new.Range("C3").Value = "Find string columns with df.select_dtypes(include=[object])"

$new.Range("A4").Value = "Convert Case"
$new.Range("B4").Value = "2 min"
$new.Range("C4").Value = "df['text_column'] = df['text_column'].str.lower() or .str.upper()"

$new.Range("C5").Value = "df['text_column'].head() to confirm the case conversion"

$new.Range("B6").Value = "10 min"

# New sheet's selection becomes the whole A1:C6 table (mirrors a
# freshly-created sheet that has just had its data filled in).
[void]$new.Range("A1:C6").Select()
